# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 19, pushing the existing
# rows 19-48 down to 20-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 19; this shifts rows 19..48 -> 20..49
$ws.Rows.Item(19).EntireRow.Insert()

# Populate the newly inserted row 19 with the new record
$ws.Cells.Item(19, 1).Value = 9
$ws.Cells.Item(19, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(19, 3).Value = "Metropolitana"
$ws.Cells.Item(19, 4).Value = Get-Date -Year 2023 -Month 5 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(19, 5).Value = 13
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100104
$ws.Cells.Item(19, 8).Value = "Frutos de pepita"
$ws.Cells.Item(19, 9).Value = 100104001
$ws.Cells.Item(19, 10).Value = "Granada"
$ws.Cells.Item(19, 11).Value = "Wonderfull"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 470
$ws.Cells.Item(19, 14).Value = 7500
$ws.Cells.Item(19, 15).Value = 8000
$ws.Cells.Item(19, 16).Value = 7734
$ws.Cells.Item(19, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(19, 18).Value = "Región Metropolitana"
$ws.Cells.Item(19, 19).Value = 516
$ws.Cells.Item(19, 20).Value = 15
